# Re-implement Min Up-/Down-Time to be cyclic
# Fill MinUpTime (col J) / MinDownTime (col K) for every generator row (7-26)
# on the "Power ThermalGen" sheet. Existing non-blank values are left as-is;
# blank cells are filled with 1 (the minimum feasible up-/down-time for a
# cyclic unit-commitment formulation). All touched cells get an integer
# ("0") number format to match the rest of the column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Power ThermalGen")

# row -> [MinUpTime, MinDownTime]
$values = @{
    7  = @(1, 1)
    8  = @(1, 1)
    9  = @(1, 1)
    10 = @(1, 1)
    11 = @(1, 1)
    12 = @(1, 1)
    13 = @(1, 1)
    14 = @(1, 1)
    15 = @(1, 1)
    16 = @(1, 1)
    17 = @(1, 1)
    18 = @(1, 1)
    19 = @(1, 1)
    20 = @(1, 1)
    21 = @(3, 2)
    22 = @(3, 1)
    23 = @(1, 4)
    24 = @(1, 1)
    25 = @(3, 2)
    26 = @(1, 1)
}

foreach ($row in 7..26) {
    $pair = $values[$row]

    $jCell = $ws.Cells.Item($row, 10)   # column J = MinUpTime
    $jCell.Value = $pair[0]
    $jCell.NumberFormat = "0"

    $kCell = $ws.Cells.Item($row, 11)   # column K = MinDownTime
    $kCell.Value = $pair[1]
    $kCell.NumberFormat = "0"
}

# Reset the lingering selection on the frozen (bottom-left) pane back to A1.
$ws.Range("A1").Select()
